$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy number formats down for the new block (A/B/C from row 1024, D from row 1023) ---
$ws.Range("A1024:C1024").Copy()
$ws.Range("A1025:C1092").PasteSpecial(-4122)
$ws.Range("D1023").Copy()
$ws.Range("D1024:D1092").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Fill in the new date / scheduled / actual values for rows 1025-1092 ---
$newRows = @(
  (44953,65,65),
  (44954,48,45),
  (44955,53,52),
  (44956,58,56),
  (44957,60,49),
  (44958,65,56),
  (44959,62,60),
  (44960,71,67),
  (44961,48,44),
  (44962,53,52),
  (44963,58,57),
  (44964,62,62),
  (44965,60,59),
  (44966,65,65),
  (44967,77,69),
  (44968,49,47),
  (44969,59,57),
  (44970,67,66),
  (44971,57,57),
  (44972,71,65),
  (44973,64,60),
  (44974,76,73),
  (44975,58,55),
  (44976,52,52),
  (44977,61,61),
  (44978,77,72),
  (44979,64,56),
  (44980,68,64),
  (44981,70,67),
  (44982,59,56),
  (44983,63,63),
  (44984,60,55),
  (44985,79,77),
  (44986,62,62),
  (44987,69,66),
  (44988,66,61),
  (44989,57,53),
  (44990,56,55),
  (44991,60,60),
  (44992,62,62),
  (44993,69,67),
  (44994,66,63),
  (44995,76,75),
  (44996,59,56),
  (44997,71,65),
  (44998,53,53),
  (44999,77,76),
  (45000,65,64),
  (45001,74,73),
  (45002,73,72),
  (45003,71,71),
  (45004,90,84),
  (45005,69,66),
  (45006,69,68),
  (45007,63,59),
  (45008,75,72),
  (45009,67,64),
  (45010,54,53),
  (45011,62,61),
  (45012,66,63),
  (45013,77,76),
  (45014,78,72),
  (45015,72,71),
  (45016,65,56),
  (45017,49,47),
  (45018,49,48),
  (45019,55,53),
  (45020,60,58)
)
$r = 1025
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- 3. Complete the missing D1024 formula, then fill the D formulas for the new rows ---
$ws.Range("D1024").Formula = "=C1024/B1024"
$ws.Range("D1025:D1088").Formula = "=C1025/B1025"
$ws.Range("D1089:D1092").Formula = "=C1089/B1089"

# --- 4. Restore view: scroll position + selection ---
$ws.Range("D1084:D1092").Select()
